$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update city names in column A (shared strings reordered) ---
$ws.Range("A1").Value = 'Datos actualizados a 18 de Abril de 2020 a las 15:52'
$ws.Range("A22").Value = 'Segovia'
$ws.Range("A23").Value = 'Leon'
$ws.Range("A24").Value = 'Asturias'
$ws.Range("A25").Value = 'Gipuzkoa/Guipuzcoa'
$ws.Range("A27").Value = 'Cantabria'
$ws.Range("A28").Value = 'A Coruña'
$ws.Range("A29").Value = 'Granada'
$ws.Range("A35").Value = 'Jaen'
$ws.Range("A36").Value = 'Cordoba'
$ws.Range("A37").Value = 'Castello/Castellon'
$ws.Range("A38").Value = 'Cuenca'
$ws.Range("A39").Value = 'Soria'
$ws.Range("A46").Value = 'Zamora'
$ws.Range("A47").Value = 'Lugo'

# --- Update numeric stats in columns B/C/D/E ---
$ws.Range("B4").Value = 52946
$ws.Range("C4").Value = 30475
$ws.Range("E4").Value = 7132
$ws.Range("E6").Value = 430
$ws.Range("E7").Value = 515
$ws.Range("B8").Value = 6212
$ws.Range("C8").Value = 3838
$ws.Range("E8").Value = 743
$ws.Range("C9").Value = 2194
$ws.Range("D9").Value = 2583
$ws.Range("E9").Value = 515
$ws.Range("B10").Value = 4579
$ws.Range("C10").Value = 954
$ws.Range("D10").Value = 3271
$ws.Range("E10").Value = 354
$ws.Range("B11").Value = 4250
$ws.Range("C11").Value = 2201
$ws.Range("E11").Value = 272
$ws.Range("B12").Value = 3831
$ws.Range("C12").Value = 3838
$ws.Range("E12").Value = 484
$ws.Range("B13").Value = 3709
$ws.Range("C13").Value = 3838
$ws.Range("E13").Value = 354
$ws.Range("C14").Value = 1677
$ws.Range("D14").Value = 1443
$ws.Range("E14").Value = 388
$ws.Range("C15").Value = 869
$ws.Range("D15").Value = 2258
$ws.Range("E15").Value = 458
$ws.Range("E16").Value = 308
$ws.Range("B17").Value = 3098
$ws.Range("C17").Value = 1045
$ws.Range("E17").Value = 243
$ws.Range("B19").Value = 2514
$ws.Range("C19").Value = 761
$ws.Range("E19").Value = 276
$ws.Range("B20").Value = 2472
$ws.Range("C20").Value = 835
$ws.Range("D20").Value = 1420
$ws.Range("E20").Value = 217
$ws.Range("B21").Value = 2299
$ws.Range("C21").Value = 432
$ws.Range("D21").Value = 1660
$ws.Range("E21").Value = 207
$ws.Range("B22").Value = 2285
$ws.Range("C22").Value = 636
$ws.Range("D22").Value = 1469
$ws.Range("E22").Value = 167
$ws.Range("B23").Value = 2285
$ws.Range("C23").Value = 1031
$ws.Range("D23").Value = 927
$ws.Range("E23").Value = 290
$ws.Range("B24").Value = 2272
$ws.Range("C24").Value = 575
$ws.Range("D24").Value = 1510
$ws.Range("E24").Value = 187
$ws.Range("B25").Value = 2266
$ws.Range("C25").Value = 6144
$ws.Range("D25").Value = 4953
$ws.Range("E25").Value = 197
$ws.Range("C26").Value = 340
$ws.Range("D26").Value = 1503
$ws.Range("E26").Value = 316
$ws.Range("B27").Value = 1990
$ws.Range("C27").Value = 534
$ws.Range("D27").Value = 1307
$ws.Range("E27").Value = 149
$ws.Range("B28").Value = 1969
$ws.Range("C28").Value = 333
$ws.Range("D28").Value = 1788
$ws.Range("E28").Value = 67
$ws.Range("B29").Value = 1969
$ws.Range("C29").Value = 563
$ws.Range("D29").Value = 1213
$ws.Range("E29").Value = 193
$ws.Range("B30").Value = 1638
$ws.Range("C30").Value = 885
$ws.Range("B32").Value = 1522
$ws.Range("C32").Value = 620
$ws.Range("E32").Value = 160
$ws.Range("B34").Value = 1345
$ws.Range("C34").Value = 3838
$ws.Range("E34").Value = 180
$ws.Range("B35").Value = 1274
$ws.Range("C35").Value = 252
$ws.Range("D35").Value = 891
$ws.Range("E35").Value = 131
$ws.Range("B36").Value = 1266
$ws.Range("C36").Value = 353
$ws.Range("D36").Value = 841
$ws.Range("E36").Value = 72
$ws.Range("B37").Value = 1257
$ws.Range("C37").Value = 435
$ws.Range("D37").Value = 739
$ws.Range("E37").Value = 139
$ws.Range("B38").Value = 1252
$ws.Range("C38").Value = 3838
$ws.Range("D38").Value = 10545
$ws.Range("E38").Value = 152
$ws.Range("B39").Value = 1231
$ws.Range("C39").Value = 293
$ws.Range("D39").Value = 802
$ws.Range("E39").Value = 94
$ws.Range("B40").Value = 1122
$ws.Range("C40").Value = 280
$ws.Range("D40").Value = 768
$ws.Range("E40").Value = 74
$ws.Range("B41").Value = 1090
$ws.Range("C41").Value = 446
$ws.Range("B42").Value = 972
$ws.Range("C42").Value = 396
$ws.Range("D42").Value = 504
$ws.Range("E42").Value = 72
$ws.Range("B45").Value = 673
$ws.Range("C45").Value = 220
$ws.Range("E45").Value = 58
$ws.Range("B46").Value = 595
$ws.Range("C46").Value = 210
$ws.Range("D46").Value = 314
$ws.Range("E46").Value = 63
$ws.Range("B47").Value = 586
$ws.Range("C47").Value = 333
$ws.Range("D47").Value = 520
$ws.Range("E47").Value = 11
$ws.Range("C48").Value = 128
$ws.Range("D48").Value = 371
$ws.Range("E48").Value = 78
$ws.Range("C49").Value = 113
$ws.Range("D49").Value = 360
$ws.Range("E49").Value = 64
$ws.Range("B51").Value = 444
$ws.Range("C51").Value = 126
$ws.Range("D51").Value = 277
$ws.Range("E51").Value = 41
$ws.Range("B52").Value = 358
$ws.Range("C52").Value = 107
$ws.Range("D52").Value = 219
$ws.Range("E52").Value = 32
$ws.Range("B54").Value = 108
$ws.Range("C54").Value = 53
$ws.Range("D54").Value = 51
